$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row before the current row 179 ("indst" / IESE row), shifting
# rows 179:293 down to 180:294.
$ws.Rows.Item(179).Insert()

# Populate the new row 179 with the new "indst" acronym entry: IESD.
$ws.Cells.Item(179, 1).Value2 = "indst"
$ws.Cells.Item(179, 2).Value2 = "IESD"
$ws.Cells.Item(179, 3).Value2 = "Industrial Equipment Survival Data"
$ws.Cells.Item(179, 4).Value2 = "Start Year Share of Industrial Equipment by Vintage, Retiring Fraction of Industrial Equipment Retiring by Age"
$ws.Cells.Item(179, 6).Value2 = "low"
$ws.Cells.Item(179, 7).Value2 = "`"Start Year Share of Industrial Equipment by Vintage`" must be updated if changing the first simulated year (Initial Time)"

# Row 179 was inserted with row-above formatting copied in (fill for column F
# came through as the "optional" shade); restore the "low" shading used by
# the rest of the indst block by copying format from a cell that already
# carries the "low" style (row 181, post-insert == old row 180).
$srcLow = $ws.Cells.Item(181, 6)
$srcLow.Copy()
$ws.Cells.Item(179, 6).PasteSpecial(-4122)

# Match the row height used for the other two-line wrapped rows in this block.
$ws.Rows.Item(179).RowHeight = 29

# Update the frozen-pane anchor and current selection to reflect the new
# layout position after the insert.
$ws.Application.ActiveWindow.ScrollRow = 171
$sheetView = $ws.Application.ActiveWindow
$ws.Range("C179").Select()
